# NFL 2024 results — add Week 9 and Week 10 sheets, plus a couple of
# score-correction tweaks on existing sheets (Week 5 overtime flag,
# Week 8 selection/active-tab cleanup as a result of Week 10 becoming
# the newly active sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Week 5: game on row 14 (Ravens @ Bengals) actually went to a
#    (second) overtime -- mark column G with " T" instead of plain "F".
# ---------------------------------------------------------------------
$week5 = $wb.Worksheets.Item("Week 5")
$week5.Range("G14").Value = " T"

# ---------------------------------------------------------------------
# 2) Add the "Week 9" worksheet at the end, after "Week 8".
# ---------------------------------------------------------------------
$week8 = $wb.Worksheets.Item("Week 8")
$week9 = $wb.Worksheets.Add($null, $week8)
$week9.Name = "Week 9"

$week9.Range("A1").Value = "Visitors team"
$week9.Range("B1").Value = "Visitors score"
$week9.Range("C1").Value = "Home team"
$week9.Range("D1").Value = "Home score"
$week9.Range("E1").Value = "Exhibition"
$week9.Range("F1").Value = "Prediction"
$week9.Range("G1").Value = "Overtime"
$week9.Range("A1:G1").Font.Color = 0

$week9Rows = @(
    @(2, "Texans", "13", "Jets", "21", "F", "Texans", "F", $true),
    @(3, "Lions", "24", "Packers", "14", "F", "Lions", "F", $false),
    @(4, "Saints", "22", "Panthers", "23", "F", "Saints", "F", $true),
    @(5, "Commanders", "27", "Giants", "22", "F", "Commanders", "F", $false),
    @(6, "Patriots", "17", "Titans", "20", "F", "Titans", "T", $false),
    @(7, "Chargers", "27", "Browns", "10", "F", "Chargers", "F", $false),
    @(8, "Raiders", "24", "Bengals", "41", "F", "Bengals", "F", $false),
    @(9, "Dolphins", "27", "Bills", "30", "F", "Bills", "F", $false),
    @(10, "Cowboys", "21", "Falcons", "27", "F", "Falcons", "F", $false),
    @(11, "Jaguars", "23", "Eagles", "28", "F", "Eagles", "F", $false),
    @(12, "Bears", "9", "Cardinals", "29", "F", "Cardinals", "F", $false),
    @(13, "Rams", "26", "Seahawks", "20", "F", "Rams", "T", $false),
    @(14, "Colts", "13", "Vikings", "21", "F", "Vikings", "F", $false),
    @(15, "Buccaneers", "24", "Chiefs", "30", "F", "Chiefs", "T", $false),
    @(16, "Broncos", "10", "Ravens", "41", "F", "Ravens", "F", $false)
)

for ($i = 0; $i -lt $week9Rows.Count; $i++) {
    $row = $week9Rows[$i]
    $r = $row[0]
    $week9.Cells.Item($r, 1).Value = $row[1]
    if ($row[2] -ne $null) { $week9.Cells.Item($r, 2).Value = [double]$row[2] }
    $week9.Cells.Item($r, 3).Value = $row[3]
    if ($row[4] -ne $null) { $week9.Cells.Item($r, 4).Value = [double]$row[4] }
    $week9.Cells.Item($r, 5).Value = $row[5]
    if ($row[6] -ne $null) {
        $week9.Cells.Item($r, 6).Value = $row[6]
        if ($row[8]) { $week9.Cells.Item($r, 6).Font.Color = 255 }
    }
    if ($row[7] -ne $null) { $week9.Cells.Item($r, 7).Value = $row[7] }
}

$week9.Range("A1:G1").Select()

# ---------------------------------------------------------------------
# 3) Add the "Week 10" worksheet at the end, after "Week 9".
# ---------------------------------------------------------------------
$week10 = $wb.Worksheets.Add($null, $week9)
$week10.Name = "Week 10"

$week10.Range("A1").Value = "Visitors team"
$week10.Range("B1").Value = "Visitors score"
$week10.Range("C1").Value = "Home team"
$week10.Range("D1").Value = "Home score"
$week10.Range("E1").Value = "Exhibition"
$week10.Range("F1").Value = "Prediction"
$week10.Range("G1").Value = "Overtime"
$week10.Range("A1:G1").Font.Color = 0

$week10Rows = @(
    @(2, "Bengals", "34", "Ravens", "35", "F", "Ravens", "F", $false),
    @(3, "Giants", "17", "Panthers", "20", "T", "Giants", "T", $true),
    @(4, "Vikings", $null, "Jaguars", $null, "F", "Vikings", $null, $false),
    @(5, "49ers", $null, "Buccaneers", $null, "F", "Buccaneers", $null, $false),
    @(6, "Falcons", $null, "Saints", $null, "F", "Falcons", $null, $false),
    @(7, "Broncos", "14", "Chiefs", "16", "F", "Chiefs", "F", $false),
    @(8, "Bills", $null, "Colts", $null, "F", $null, $null, $false),
    @(9, "Patriots", $null, "Bears", $null, "F", "Bears", $null, $false),
    @(10, "Titans", $null, "Chargers", $null, "F", "Chargers", $null, $false),
    @(11, "Jets", $null, "Cardinals", $null, "F", $null, $null, $false),
    @(12, "Lions", $null, "Texans", $null, "F", "Lions", $null, $false),
    @(13, "Dolphins", $null, "Rams", $null, "F", "Rams", $null, $false),
    @(14, "Steelers", "28", "Commanders", "27", "F", "Commanders", "F", $true)
)

for ($i = 0; $i -lt $week10Rows.Count; $i++) {
    $row = $week10Rows[$i]
    $r = $row[0]
    $week10.Cells.Item($r, 1).Value = $row[1]
    if ($row[2] -ne $null) { $week10.Cells.Item($r, 2).Value = [double]$row[2] }
    $week10.Cells.Item($r, 3).Value = $row[3]
    if ($row[4] -ne $null) { $week10.Cells.Item($r, 4).Value = [double]$row[4] }
    $week10.Cells.Item($r, 5).Value = $row[5]
    if ($row[6] -ne $null) {
        $week10.Cells.Item($r, 6).Value = $row[6]
        if ($row[8]) { $week10.Cells.Item($r, 6).Font.Color = 255 }
    }
    if ($row[7] -ne $null) { $week10.Cells.Item($r, 7).Value = $row[7] }
}

$week10.Range("F3").Select()
$week10.Activate()

# ---------------------------------------------------------------------
# 4) Workbook window: the newly-added "Week 10" tab becomes the active
#    one (this falls out naturally from activating $week10 above, but
#    we set it explicitly too for clarity).
# ---------------------------------------------------------------------
$excel.ActiveWindow.ActiveSheet
